# Update 'F' column (想去人数 / want-to-go count) values across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 8118
$ws.Cells.Item(8, 6).Value = 2102
$ws.Cells.Item(9, 6).Value = 66
$ws.Cells.Item(12, 6).Value = 504
$ws.Cells.Item(13, 6).Value = 1063
$ws.Cells.Item(15, 6).Value = 149
$ws.Cells.Item(16, 6).Value = 1148
$ws.Cells.Item(18, 6).Value = 722
$ws.Cells.Item(19, 6).Value = 508
$ws.Cells.Item(20, 6).Value = 5
$ws.Cells.Item(22, 6).Value = 416
$ws.Cells.Item(23, 6).Value = 4703
$ws.Cells.Item(24, 6).Value = 105
$ws.Cells.Item(25, 6).Value = 49892
$ws.Cells.Item(26, 6).Value = 4025
$ws.Cells.Item(28, 6).Value = 990
$ws.Cells.Item(29, 6).Value = 751
$ws.Cells.Item(30, 6).Value = 332
$ws.Cells.Item(31, 6).Value = 67
$ws.Cells.Item(32, 6).Value = 832
$ws.Cells.Item(35, 6).Value = 192
$ws.Cells.Item(40, 6).Value = 1010
$ws.Cells.Item(43, 6).Value = 1048
$ws.Cells.Item(44, 6).Value = 677
$ws.Cells.Item(48, 6).Value = 23
$ws.Cells.Item(49, 6).Value = 2452

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 241
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(12, 6).Value = 109
$ws.Cells.Item(14, 6).Value = 44
$ws.Cells.Item(16, 6).Value = 82
$ws.Cells.Item(18, 6).Value = 27
$ws.Cells.Item(19, 6).Value = 150
$ws.Cells.Item(20, 6).Value = 7301
$ws.Cells.Item(21, 6).Value = 64

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 2203
$ws.Cells.Item(5, 6).Value = 1482
$ws.Cells.Item(8, 6).Value = 2321
$ws.Cells.Item(9, 6).Value = 9242
$ws.Cells.Item(10, 6).Value = 1520
$ws.Cells.Item(12, 6).Value = 60
$ws.Cells.Item(14, 6).Value = 41

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 2203
$ws.Cells.Item(4, 6).Value = 8118
$ws.Cells.Item(5, 6).Value = 1482
$ws.Cells.Item(7, 6).Value = 1520
$ws.Cells.Item(10, 6).Value = 66
$ws.Cells.Item(13, 6).Value = 1063
$ws.Cells.Item(14, 6).Value = 241
$ws.Cells.Item(18, 6).Value = 149
$ws.Cells.Item(19, 6).Value = 1148
$ws.Cells.Item(21, 6).Value = 416
$ws.Cells.Item(22, 6).Value = 4703
$ws.Cells.Item(23, 6).Value = 105
$ws.Cells.Item(24, 6).Value = 109
$ws.Cells.Item(25, 6).Value = 44
$ws.Cells.Item(26, 6).Value = 4025
$ws.Cells.Item(28, 6).Value = 990
$ws.Cells.Item(29, 6).Value = 751
$ws.Cells.Item(30, 6).Value = 332
$ws.Cells.Item(31, 6).Value = 67
$ws.Cells.Item(32, 6).Value = 832
$ws.Cells.Item(35, 6).Value = 192
$ws.Cells.Item(41, 6).Value = 1048
$ws.Cells.Item(42, 6).Value = 677
$ws.Cells.Item(47, 6).Value = 23
